$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab to match the new table title.
$ws.Name = "RT-qPCR_dCT"

# "gene TSS" -> "gene TSS-Mu" everywhere it appears in column C.
for ($r = 1; $r -le 35; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Text -eq "gene TSS") {
        $cell.Value = "gene TSS-Mu"
    }
}

# Column C needs to be a bit wider to fit the longer "gene TSS-Mu" label
# (closest attainable width to the recorded best-fit result of 13.6640625).
$ws.Columns.Item(3).ColumnWidth = 12.8

# Move the selection/scroll position that was saved with the workbook.
$null = $ws.Range("E16").Select()
